$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$startRow = $used.Row

for ($i = 0; $i -lt $rowCount; $i++) {
    $r = $startRow + $i
    $cell = $ws.Cells.Item($r, 7)  # Column G
    if ($cell.Text -eq "System, dnasr281@gmail.com") {
        $cell.Value = "dnasr281@gmail.com, System"
    }
}
